$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.370.45"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.906.99"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.58%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0720"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "2.179.15"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.699"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.903.52"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("D18").Value = "35.377.71"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "240.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.06%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0566"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.934"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.67%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0634"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "90.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "1.340.96"
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +38.67%  "
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.03%  "
$ws.Range("D51").Value = "2.092.83"
$ws.Range("E51").Value = "  +2.46%  "
